$d = $word.ActiveDocument
$d.Content.Find.Execute("tool.  ", $true, $false, $false, $false, $false, $true, 1, $false, "tool, and the complete methodology document can be downloaded using the button below. ", 2)
